$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '34.545.07'
$ws.Range('E2').Value = '  +0.45%  '

# Row 3
$ws.Range('D3').Value = '1.804.07'
$ws.Range('E3').Value = '  +0.18%  '

# Row 4
$ws.Range('E4').Value = '  +0.19%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '224.81'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.26%  '

# Row 6
$ws.Range('E6').Value = '  -0.10%  '

# Row 7
$ws.Range('E7').Value = '  +0.11%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '42.12'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +16.74%  '

# Row 9
$ws.Range('E9').Value = '  +0.27%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0667'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -1.83%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0996'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +3.20%  '

# Row 12
$ws.Range('D12').Value = '2.062.18'
$ws.Range('E12').Value = '  +0.06%  '

# Row 13
$ws.Range('D13').Value = '1.796.30'
$ws.Range('E13').Value = '  -0.48%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '10.90'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -2.85%  '

# Row 15
$ws.Range('D15').Value = '34.517.72'
$ws.Range('E15').Value = '  +0.42%  '

# Row 16
$ws.Range('E16').Value = '  -0.31%  '

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '4.40'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.63%  '

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '67.28'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -1.83%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '240.26'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -1.10%  '

# Row 20
$ws.Range('D20').Value = '0.0₃0769'
$ws.Range('E20').Value = '  -0.63%  '

# Row 21
$ws.Range('E21').Value = '  -0.98%  '

# Row 22
$ws.Range('E22').Value = '  +0.21%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.36'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +6.69%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.14'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -2.79%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '170.71'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.38%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.67'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -2.72%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '17.45'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.63%  '

# Row 28
$ws.Range('E28').Value = '  -0.37%  '

# Row 29
$ws.Range('E29').Value = '  +0.15%  '

# Row 30
$ws.Range('E30').Value = '  +0.13%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.23'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.17%  '

# Row 32
$ws.Range('E32').Value = '  -1.01%  '

# Row 33
$ws.Range('E33').Value = '  -0.64%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.79'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.81%  '

# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '87.77'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +8.38%  '

# Row 36
$ws.Range('E36').Value = '  -0.20%  '

# Row 37
$ws.Range('D37').Value = '1.316.71'
$ws.Range('E37').Value = '  -3.45%  '

# Row 38
$ws.Range('E38').Value = '  +0.34%  '

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '14.93'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +13.22%  '

# Row 40
$ws.Range('E40').Value = '  +0.74%  '

# Row 41
$ws.Range('E41').Value = '  -1.65%  '

# Row 42
$ws.Range('E42').Value = '  +4.94%  '

# Row 43
$ws.Range('E43').Value = '  +0.19%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.80'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.42%  '

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.937'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.09%  '

# Row 46
$ws.Range('E46').Value = '  +4.07%  '

# Row 47
$ws.Range('D47').Value = '1.962.94'
$ws.Range('E47').Value = '  +0.03%  '

# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '5.80'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.21%  '

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '100.44'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -1.48%  '

# Row 51
$ws.Range('E51').Value = '  +0.71%  '
